# Update cryptocurrency price and volume data in the active worksheet.
# Force Price (D) and Volume(1h) (E) columns to Text format first so that
# numeric-looking strings (e.g. "1.000", "244.60") are preserved exactly
# as text instead of being auto-converted/normalized to numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '30.383.18'
$ws.Range("E2").Value = '  -1.42%  '
$ws.Range("D3").Value = '1.871.49'
$ws.Range("E3").Value = '  -1.11%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '244.60'
$ws.Range("E5").Value = '  -1.98%  '
$ws.Range("D6").Value = '1.000'
$ws.Range("E7").Value = '  -1.08%  '
$ws.Range("D8").Value = '0.2873'
$ws.Range("E8").Value = '  -2.35%  '
$ws.Range("D9").Value = '0.06487'
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("D10").Value = '21.82'
$ws.Range("E10").Value = '  -1.28%  '
$ws.Range("D11").Value = '100.12'
$ws.Range("E11").Value = '  +2.71%  '
$ws.Range("E12").Value = '  +0.27%  '
$ws.Range("D13").Value = '1.871.63'
$ws.Range("E13").Value = '  -1.02%  '
$ws.Range("D14").Value = '0.7271'
$ws.Range("E14").Value = '  -1.60%  '
$ws.Range("D15").Value = '5.171'
$ws.Range("E15").Value = '  -1.57%  '
$ws.Range("D16").Value = '285.14'
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("D17").Value = '30.377.53'
$ws.Range("E17").Value = '  -1.49%  '
$ws.Range("D18").Value = '13.08'
$ws.Range("E18").Value = '  -1.07%  '
$ws.Range("D19").Value = '0.9999'
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("D20").Value = '0.000007493'
$ws.Range("E20").Value = '  -1.37%  '
$ws.Range("D21").Value = '2.115.97'
$ws.Range("E21").Value = '  -0.97%  '
$ws.Range("D22").Value = '5.340'
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").Value = '0.9996'
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '6.315'
$ws.Range("E24").Value = '  +0.77%  '
$ws.Range("D25").Value = '163.14'
$ws.Range("E25").Value = '  -0.82%  '
$ws.Range("D26").Value = '9.032'
$ws.Range("E26").Value = '  -2.57%  '
$ws.Range("D27").Value = '18.97'
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("D28").Value = '1.896'
$ws.Range("E28").Value = '  -1.95%  '
$ws.Range("D29").Value = '0.09661'
$ws.Range("E29").Value = '  -0.91%  '
$ws.Range("D30").Value = '1.320'
$ws.Range("E30").Value = '  -1.88%  '
$ws.Range("D31").Value = '1.487'
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("D32").Value = '4.223'
$ws.Range("E32").Value = '  -2.30%  '
$ws.Range("D33").Value = '4.140'
$ws.Range("E33").Value = '  -1.45%  '
$ws.Range("D34").Value = '0.04806'
$ws.Range("E34").Value = '  -1.73%  '
$ws.Range("D35").Value = '1.125'
$ws.Range("E35").Value = '  -0.56%  '
$ws.Range("D36").Value = '0.6883'
$ws.Range("D37").Value = '2.725'
$ws.Range("E37").Value = '  +0.18%  '
$ws.Range("D38").Value = '0.01902'
$ws.Range("E38").Value = '  -0.71%  '
$ws.Range("D39").Value = '2.841'
$ws.Range("E39").Value = '  +1.21%  '
$ws.Range("D40").Value = '76.01'
$ws.Range("E40").Value = '  -0.04%  '
$ws.Range("D41").Value = '6.276'
$ws.Range("E41").Value = '  -1.43%  '
$ws.Range("D42").Value = '1.960'
$ws.Range("E42").Value = '  -3.70%  '
$ws.Range("D43").Value = '0.4217'
$ws.Range("E43").Value = '  -1.63%  '
$ws.Range("D44").Value = '0.9991'
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").Value = '0.8240'
$ws.Range("E45").Value = '  -1.89%  '
$ws.Range("D46").Value = '101.08'
$ws.Range("E46").Value = '  -1.00%  '
$ws.Range("D47").Value = '9.755'
$ws.Range("E47").Value = '  +3.52%  '
$ws.Range("D48").Value = '7.012'
$ws.Range("E48").Value = '  -1.39%  '
$ws.Range("D49").Value = '35.10'
$ws.Range("E49").Value = '  -2.15%  '
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("D51").Value = '890.27'
$ws.Range("E51").Value = '  -4.14%  '
